# chore: update Sheets via scheduled runner
# Refresh cached market-board price snapshots (currentAveragePrice* / Leve profit columns)
# across the per-job Leve worksheets (Goblin_Profits workbook).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 332.77777
$ws.Range("I2").Value = 332.77777
$ws.Range("K2").Value = 332.77777
$ws.Range("M2").Value = -219.77777

$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2770

$ws.Range("H70").Value = 45837748
$ws.Range("I70").Value = 12502996
$ws.Range("J70").Value = 62505124
$ws.Range("K70").Value = 37508988
$ws.Range("L70").Value = 187515372
$ws.Range("M70").Value = -37508718
$ws.Range("N70").Value = -187515912

$ws.Range("H73").Value = 45837748
$ws.Range("I73").Value = 12502996
$ws.Range("J73").Value = 62505124
$ws.Range("K73").Value = 37508988
$ws.Range("L73").Value = 187515372
$ws.Range("M73").Value = -37508052
$ws.Range("N73").Value = -187517244

$ws.Range("H95").Value = 47000
$ws.Range("J95").Value = 47000
$ws.Range("L95").Value = 47000
$ws.Range("N95").Value = -52492

$ws.Range("H98").Value = 9570.954
$ws.Range("I98").Value = 12572.267
$ws.Range("J98").Value = 3139.5715
$ws.Range("K98").Value = 12572.267
$ws.Range("L98").Value = 3139.5715
$ws.Range("M98").Value = -11074.267
$ws.Range("N98").Value = -6135.5715

$ws.Range("H122").Value = 9570.954
$ws.Range("I122").Value = 12572.267
$ws.Range("J122").Value = 3139.5715
$ws.Range("K122").Value = 37716.801
$ws.Range("L122").Value = 9418.7145
$ws.Range("M122").Value = -35266.801
$ws.Range("N122").Value = -14318.7145

$ws.Range("H137").Value = 4035.5588
$ws.Range("I137").Value = 4535.8213
$ws.Range("K137").Value = 13607.4639
$ws.Range("M137").Value = -11057.4639

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2143.4312
$ws.Range("I32").Value = 2093.3157
$ws.Range("K32").Value = 2093.3157
$ws.Range("M32").Value = -1806.3157

$ws.Range("H50").Value = 7889.9
$ws.Range("I50").Value = 416
$ws.Range("K50").Value = 416
$ws.Range("M50").Value = 298

$ws.Range("H76").Value = 46000
$ws.Range("J76").Value = 46000
$ws.Range("L76").Value = 46000
$ws.Range("N76").Value = -46676

$ws.Range("H79").Value = 46000
$ws.Range("J79").Value = 46000
$ws.Range("L79").Value = 46000
$ws.Range("N79").Value = -48340

$ws.Range("H95").Value = 36241.6
$ws.Range("J95").Value = 36241.6
$ws.Range("L95").Value = 36241.6
$ws.Range("N95").Value = -41733.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3400.2307
$ws.Range("I2").Value = 2587.25
$ws.Range("J2").Value = 4701
$ws.Range("K2").Value = 2587.25
$ws.Range("L2").Value = 4701
$ws.Range("M2").Value = -2474.25
$ws.Range("N2").Value = -4927

$ws.Range("H22").Value = 1174.4839
$ws.Range("I22").Value = 896.8823
$ws.Range("J22").Value = 1511.5714
$ws.Range("K22").Value = 896.8823
$ws.Range("L22").Value = 1511.5714
$ws.Range("M22").Value = -546.8823
$ws.Range("N22").Value = -2211.5714

$ws.Range("H50").Value = 58984.5
$ws.Range("J50").Value = 58984.5
$ws.Range("L50").Value = 58984.5
$ws.Range("N50").Value = -60234.5

$ws.Range("H99").Value = 2883.3
$ws.Range("I99").Value = 2943.6667
$ws.Range("J99").Value = 2792.75
$ws.Range("K99").Value = 2943.6667
$ws.Range("L99").Value = 2792.75
$ws.Range("M99").Value = -1445.6667
$ws.Range("N99").Value = -5788.75

$ws.Range("H126").Value = 2883.3
$ws.Range("I126").Value = 2943.6667
$ws.Range("J126").Value = 2792.75
$ws.Range("K126").Value = 8831.000100000001
$ws.Range("L126").Value = 8378.25
$ws.Range("M126").Value = -6361.000100000001
$ws.Range("N126").Value = -13318.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 4498.222
$ws.Range("I13").Value = 314.66666
$ws.Range("J13").Value = 12865.333
$ws.Range("K13").Value = 943.9999799999999
$ws.Range("L13").Value = 38595.999
$ws.Range("M13").Value = -775.9999799999999
$ws.Range("N13").Value = -38931.999

$ws.Range("H22").Value = 366.66666
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 1200
$ws.Range("M22").Value = -1031

$ws.Range("H27").Value = 366.66666
$ws.Range("I27").Value = 400
$ws.Range("K27").Value = 1200
$ws.Range("M27").Value = -1098

$ws.Range("H49").Value = 1080
$ws.Range("I49").Value = 1201.5
$ws.Range("K49").Value = 3604.5
$ws.Range("M49").Value = -3448.5

$ws.Range("H60").Value = 222.85715
$ws.Range("I60").Value = 169.5
$ws.Range("J60").Value = 543
$ws.Range("K60").Value = 508.5
$ws.Range("L60").Value = 1629
$ws.Range("M60").Value = -257.5
$ws.Range("N60").Value = -2131

$ws.Range("H131").Value = 2022290.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 31187.5
$ws.Range("J26").Value = 31187.5
$ws.Range("L26").Value = 31187.5
$ws.Range("N26").Value = -31747.5

$ws.Range("H50").Value = 31187.5
$ws.Range("J50").Value = 31187.5
$ws.Range("L50").Value = 31187.5
$ws.Range("N50").Value = -32183.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 40000
$ws.Range("J96").Value = 40000
$ws.Range("L96").Value = 40000
$ws.Range("N96").Value = -45492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39998.5
$ws.Range("J70").Value = 39998.5
$ws.Range("L70").Value = 39998.5
$ws.Range("N70").Value = -40628.5

$ws.Range("H73").Value = 39998.5
$ws.Range("J73").Value = 39998.5
$ws.Range("L73").Value = 39998.5
$ws.Range("N73").Value = -42182.5

$ws.Range("H81").Value = 1998.8572
$ws.Range("I81").Value = 1998.4
$ws.Range("K81").Value = 3996.8
$ws.Range("M81").Value = -2935.8

$ws.Range("H84").Value = 1998.8572
$ws.Range("I84").Value = 1998.4
$ws.Range("K84").Value = 19984
$ws.Range("M84").Value = -14680

$ws.Range("H99").Value = 33499.5
$ws.Range("I99").Value = 16999
$ws.Range("K99").Value = 16999
$ws.Range("M99").Value = -14004

$ws.Range("H113").Value = 2397.4
$ws.Range("J113").Value = 2750
$ws.Range("L113").Value = 8250
$ws.Range("N113").Value = -12590

$ws.Range("H126").Value = 4531.615
$ws.Range("I126").Value = 5268.5
$ws.Range("J126").Value = 3352.6
$ws.Range("K126").Value = 15805.5
$ws.Range("L126").Value = 10057.8
$ws.Range("M126").Value = -13335.5
$ws.Range("N126").Value = -14997.8
